$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) values that look numeric (single decimal point) must be forced to
# Text format first, otherwise Excel auto-converts them to numbers and mangles
# the display (e.g. floating point artifacts like 313.41000000000003).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.529.00"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "2.479.48"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "313.41"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "92.86"
$ws.Range("E6").Value = "  -1.43%  "

$ws.Range("E7").Value = "  -1.21%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  +1.86%  "

$ws.Range("D10").Value = "32.74"
$ws.Range("E10").Value = "  -1.85%  "

$ws.Range("E12").Value = "  +2.04%  "

$ws.Range("D13").Value = "2.861.28"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("D15").Value = "16.12"
$ws.Range("E15").Value = "  +8.56%  "

$ws.Range("D16").Value = "2.503.32"
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("E17").Value = "  -2.12%  "

$ws.Range("D18").Value = "41.522.74"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").Value = "71.72"
$ws.Range("E21").Value = "  +5.03%  "

$ws.Range("D22").Value = "11.26"
$ws.Range("E22").Value = "  -0.58%  "

$ws.Range("D23").Value = "236.59"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D27").Value = "24.91"
$ws.Range("E27").Value = "  +3.83%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").Value = "35.99"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("D34").Value = "0.0756"
$ws.Range("E34").Value = "  +1.67%  "

$ws.Range("E35").Value = "  -8.30%  "

$ws.Range("D36").Value = "17.22"
$ws.Range("E36").Value = "  +1.41%  "

$ws.Range("E37").Value = "  +3.73%  "

$ws.Range("E38").Value = "  -4.60%  "

$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -3.32%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").Value = "19.49"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("D44").Value = "1.979.73"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("E47").Value = "  +3.23%  "

$ws.Range("D48").Value = "2.719.71"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").Value = "97.52"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Value = "68.02"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("D51").Value = "72.28"
$ws.Range("E51").Value = "  -2.05%  "
